$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.024324031222992
$ws.Range("D2").Value = 1.028723914620983
$ws.Range("E2").Value = 1.04807214688726
$ws.Range("F2").Value = 1.052258118722738
$ws.Range("I2").Value = 1.031797358698314
$ws.Range("J2").Value = 1.029499284313648
$ws.Range("K2").Value = 1.031539900414579
$ws.Range("L2").Value = 1.05083301330704
$ws.Range("M2").Value = 1.05500734930304
$ws.Range("N2").Value = 1.013894279400376
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.025054188332102
$ws.Range("D3").Value = 1.029240845474762
$ws.Range("E3").Value = 1.049246818897453
$ws.Range("F3").Value = 1.053479227526721
$ws.Range("I3").Value = 1.031913723182578
$ws.Range("J3").Value = 1.02986943401778
$ws.Range("K3").Value = 1.031865611888017
$ws.Range("L3").Value = 1.051818718007142
$ws.Range("M3").Value = 1.056040219859817
$ws.Range("N3").Value = 1.014016929603322
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.025527031369525
$ws.Range("D4").Value = 1.029575631742335
$ws.Range("E4").Value = 1.050007975719168
$ws.Range("F4").Value = 1.054270459170101
$ws.Range("I4").Value = 1.031987993151654
$ws.Range("J4").Value = 1.030108637225587
$ws.Range("K4").Value = 1.032075952316551
$ws.Range("L4").Value = 1.05245702430219
$ws.Range("M4").Value = 1.056709081174373
$ws.Range("N4").Value = 1.014096178127843
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.025725904572258
$ws.Range("D5").Value = 1.029716445477088
$ws.Range("E5").Value = 1.050328221365174
$ws.Range("F5").Value = 1.054603354712269
$ws.Range("I5").Value = 1.032018970137097
$ws.Range("J5").Value = 1.030209123449312
$ws.Range("K5").Value = 1.032164278657341
$ws.Range("L5").Value = 1.052725485600443
$ws.Range("M5").Value = 1.056990396220991
$ws.Range("N5").Value = 1.014129466484128
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.025759301488703
$ws.Range("D6").Value = 1.029740092756243
$ws.Range("E6").Value = 1.050382006997433
$ws.Range("F6").Value = 1.05465926469404
$ws.Range("I6").Value = 1.032024156860694
$ws.Range("J6").Value = 1.030225991113546
$ws.Range("K6").Value = 1.032179103099336
$ws.Range("L6").Value = 1.052770568318891
$ws.Range("M6").Value = 1.05703763763126
$ws.Range("N6").Value = 1.014135054109844
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.025529688371396
$ws.Range("D7").Value = 1.029577513030962
$ws.Range("E7").Value = 1.050012253853619
$ws.Range("F7").Value = 1.054274906307836
$ws.Range("I7").Value = 1.031988408035403
$ws.Range("J7").Value = 1.030109980222579
$ws.Range("K7").Value = 1.032077132934703
$ws.Range("L7").Value = 1.05246061103395
$ws.Range("M7").Value = 1.056712839625588
$ws.Range("N7").Value = 1.014096623037858
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.024570711064403
$ws.Range("D8").Value = 1.028898551264187
$ws.Range("E8").Value = 1.048468912319123
$ws.Range("F8").Value = 1.052670572043259
$ws.Range("I8").Value = 1.031836896641323
$ws.Range("J8").Value = 1.029624441411301
$ws.Range("K8").Value = 1.031650061792839
$ws.Range("L8").Value = 1.051166035576403
$ws.Range("M8").Value = 1.055356304075712
$ws.Range("N8").Value = 1.013935753081569
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.022883871048742
$ws.Range("D9").Value = 1.027704483352559
$ws.Range("E9").Value = 1.045757496920823
$ws.Range("F9").Value = 1.049851889481349
$ws.Range("I9").Value = 1.031562084331042
$ws.Range("J9").Value = 1.028766540933972
$ws.Range("K9").Value = 1.030894361780613
$ws.Range("L9").Value = 1.048888569615985
$ws.Range("M9").Value = 1.052969932998265
$ws.Range("N9").Value = 1.013651419073943
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.02176142324865
$ws.Range("D10").Value = 1.026910114079195
$ws.Range("E10").Value = 1.04395534856131
$ws.Range("F10").Value = 1.047978373638712
$ws.Range("I10").Value = 1.031373645739896
$ws.Range("J10").Value = 1.028193105741059
$ws.Range("K10").Value = 1.030388507248725
$ws.Range("L10").Value = 1.047372762347004
$ws.Range("M10").Value = 1.0513817138064
$ws.Range("N10").Value = 1.013461304952356
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.021275911857088
$ws.Range("D11").Value = 1.026566561280796
$ws.Range("E11").Value = 1.043176289284176
$ws.Range("F11").Value = 1.047168446189507
$ws.Range("I11").Value = 1.031290816303789
$ws.Range("J11").Value = 1.027944457910599
$ws.Range("K11").Value = 1.030168992113089
$ws.Range("L11").Value = 1.046716991437158
$ws.Range("M11").Value = 1.050694634328351
$ws.Range("N11").Value = 1.013378855359468
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.021095650413464
$ws.Range("D12").Value = 1.026439014371502
$ws.Range("E12").Value = 1.042887104227913
$ws.Range("F12").Value = 1.046867800536021
$ws.Range("I12").Value = 1.031259864883771
$ws.Range("J12").Value = 1.02785204792177
$ws.Range("K12").Value = 1.030087383805232
$ws.Range("L12").Value = 1.046473496230605
$ws.Range("M12").Value = 1.05043951677157
$ws.Range("N12").Value = 1.01334821088246
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.02113431349979
$ws.Range("D13").Value = 1.026466370698668
$ws.Range("E13").Value = 1.042949126712308
$ws.Range("F13").Value = 1.046932281120647
$ws.Range("I13").Value = 1.031266512425659
$ws.Range("J13").Value = 1.027871872475657
$ws.Range("K13").Value = 1.030104892243321
$ws.Range("L13").Value = 1.046525722834107
$ws.Range("M13").Value = 1.050494236088291
$ws.Range("N13").Value = 1.013354785082335
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.021261009766997
$ws.Range("D14").Value = 1.026556016900508
$ws.Range("E14").Value = 1.043152381238224
$ws.Range("F14").Value = 1.047143590700341
$ws.Range("I14").Value = 1.031288261618013
$ws.Range("J14").Value = 1.027936820309247
$ws.Range("K14").Value = 1.030162247779493
$ws.Range("L14").Value = 1.046696862268534
$ws.Range("M14").Value = 1.050673544308351
$ws.Range("N14").Value = 1.013376322663611
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.021339082048538
$ws.Range("D15").Value = 1.026611259400258
$ws.Range("E15").Value = 1.043277638520034
$ws.Range("F15").Value = 1.047273811659197
$ws.Range("I15").Value = 1.031301637531602
$ws.Range("J15").Value = 1.027976830068733
$ws.Range("K15").Value = 1.030197577091895
$ws.Range("L15").Value = 1.046802318480339
$ws.Range("M15").Value = 1.05078403451507
$ws.Range("N15").Value = 1.013389590167206
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.021793656046527
$ws.Range("D16").Value = 1.026932923417322
$ws.Range("E16").Value = 1.044007079146227
$ws.Range("F16").Value = 1.048032153571571
$ws.Range("I16").Value = 1.031379116899437
$ws.Range("J16").Value = 1.028209600458041
$ws.Range("K16").Value = 1.030403065787033
$ws.Range("L16").Value = 1.047416295966761
$ws.Range("M16").Value = 1.051427326255653
$ws.Range("N16").Value = 1.013466774173246
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.022078937556147
$ws.Range("D17").Value = 1.027134806851738
$ws.Range("E17").Value = 1.044464981089511
$ws.Range("F17").Value = 1.048508193774915
$ws.Range("I17").Value = 1.031427387600766
$ws.Range("J17").Value = 1.028355519055145
$ws.Range("K17").Value = 1.030531836487604
$ws.Range("L17").Value = 1.047801583811025
$ws.Range("M17").Value = 1.051831014611914
$ws.Range("N17").Value = 1.013515155386319
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.022245387105351
$ws.Range("D18").Value = 1.027252601911002
$ws.Range("E18").Value = 1.044732191355068
$ws.Range("F18").Value = 1.048785987098581
$ws.Range("I18").Value = 1.031455423956373
$ws.Range("J18").Value = 1.028440597397128
$ws.Range("K18").Value = 1.030606900187532
$ws.Range("L18").Value = 1.04802637202976
$ws.Range("M18").Value = 1.052066539996269
$ws.Range("N18").Value = 1.013543362867187
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.022302150476768
$ws.Range("D19").Value = 1.027292773694154
$ws.Range("E19").Value = 1.044823324165427
$ws.Range("F19").Value = 1.048880729019667
$ws.Range("I19").Value = 1.031464963406015
$ws.Range("J19").Value = 1.028469601200523
$ws.Range("K19").Value = 1.030632487123901
$ws.Range("L19").Value = 1.048103028570772
$ws.Range("M19").Value = 1.052146858386625
$ws.Range("N19").Value = 1.013552978764994
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.022048324400797
$ws.Range("D20").Value = 1.027113142534248
$ws.Range("E20").Value = 1.04441583974178
$ws.Range("F20").Value = 1.048457106020023
$ws.Range("I20").Value = 1.031422220929266
$ws.Range("J20").Value = 1.028339866829477
$ws.Range("K20").Value = 1.030518025368094
$ws.Range("L20").Value = 1.047760240247011
$ws.Range("M20").Value = 1.051787696385282
$ws.Range("N20").Value = 1.013509965826579
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.02122369867765
$ws.Range("D21").Value = 1.026529616551379
$ws.Range("E21").Value = 1.043092522543246
$ws.Range("F21").Value = 1.047081359848432
$ws.Range("I21").Value = 1.031281862122047
$ws.Range("J21").Value = 1.027917696194316
$ws.Range("K21").Value = 1.030145359944719
$ws.Range("L21").Value = 1.046646463563511
$ws.Range("M21").Value = 1.0506207399043
$ws.Range("N21").Value = 1.013369980905234
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.020705682219873
$ws.Range("D22").Value = 1.026163101103676
$ws.Range("E22").Value = 1.042261612544317
$ws.Range("F22").Value = 1.046217515875662
$ws.Range("I22").Value = 1.03119254359003
$ws.Range("J22").Value = 1.027651965507306
$ws.Range("K22").Value = 1.029910642362426
$ws.Range("L22").Value = 1.045946692910424
$ws.Range("M22").Value = 1.049887573814538
$ws.Range("N22").Value = 1.013281856912888
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.020980248556006
$ws.Range("D23").Value = 1.026357362209856
$ws.Range("E23").Value = 1.042701988407584
$ws.Range("F23").Value = 1.046675347783588
$ws.Range("I23").Value = 1.031239994212244
$ws.Range("J23").Value = 1.027792862038493
$ws.Range("K23").Value = 1.030035108931498
$ws.Range("L23").Value = 1.04631760683408
$ws.Range("M23").Value = 1.050276187391023
$ws.Range("N23").Value = 1.013328583413317
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.022062157014685
$ws.Range("D24").Value = 1.027122931583818
$ws.Range("E24").Value = 1.044438044218784
$ws.Range("F24").Value = 1.048480189984391
$ws.Range("I24").Value = 1.031424555893871
$ws.Range("J24").Value = 1.028346939500154
$ws.Range("K24").Value = 1.030524266160982
$ws.Range("L24").Value = 1.04777892144971
$ws.Range("M24").Value = 1.051807269839254
$ws.Range("N24").Value = 1.01351231080379
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.023319594275587
$ws.Range("D25").Value = 1.028012890157709
$ws.Range("E25").Value = 1.046457498850742
$ws.Range("F25").Value = 1.050579597383002
$ws.Range("I25").Value = 1.031634054282378
$ws.Range("J25").Value = 1.028988597645559
$ws.Range("K25").Value = 1.031090094911725
$ws.Range("L25").Value = 1.049476907060314
$ws.Range("M25").Value = 1.053586391684787
$ws.Range("N25").Value = 1.013725026100421
